$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08658538737642299
$ws.Range("H2").Value = -34.11588900330445
$ws.Range("I2").Value = 11.09720476824272

$ws.Range("G3").Value = 0.1003643442651757
$ws.Range("H3").Value = 12.72181564629117

$ws.Range("G4").Value = -0.6362626895486008
$ws.Range("H4").Value = -3.39960620511847

$ws.Range("G5").Value = -0.635638224653324
$ws.Range("H5").Value = -4.101745327378339

$ws.Range("G6").Value = 0.1713750758374028
$ws.Range("H6").Value = -30.36584484041292

$ws.Range("G7").Value = 0.3806497397478641
$ws.Range("H7").Value = 132.3557719809244

$ws.Range("G8").Value = 0.1118397479653579
$ws.Range("H8").Value = -32.34500455713895

$ws.Range("G9").Value = 0.2124975058568356
$ws.Range("H9").Value = 8.921390360598384

$ws.Range("G10").Value = -0.1360603570835038
$ws.Range("H10").Value = -138.0910814425618

$ws.Range("G11").Value = -0.09776905627177691
$ws.Range("H11").Value = 17.68086490856194

$ws.Range("G12").Value = 0.1747859421836568
$ws.Range("H12").Value = 9.908021853321596

$ws.Range("G13").Value = 0.2460904436293652
$ws.Range("H13").Value = 19.65785901368248

$ws.Range("G14").Value = 0.2087152063909201
$ws.Range("H14").Value = 10.21376980133935

$ws.Range("G15").Value = 0.2301319581552699
$ws.Range("H15").Value = -7.912377271219825

$ws.Range("G16").Value = 0.05385407473209116
$ws.Range("H16").Value = 47.6253022292823

$ws.Range("G17").Value = 0.0249146863594033
$ws.Range("H17").Value = -29.75959428415728

$ws.Range("G18").Value = 0.1426927472281291
$ws.Range("H18").Value = -17.66818415287008

$ws.Range("G19").Value = 0.1659454741328804
$ws.Range("H19").Value = 31.96150248796986

$ws.Range("G20").Value = 0.1167005661561528
$ws.Range("H20").Value = 1.786768076873927

$ws.Range("G21").Value = 0.1602754673295104
$ws.Range("H21").Value = 59.63878643684927

$ws.Range("G22").Value = 0.08632827573073851
$ws.Range("H22").Value = -8.353602164188718

$ws.Range("G23").Value = 0.07187625681153532
$ws.Range("H23").Value = -33.74848847853746

$ws.Range("G24").Value = -0.1991847205090597
$ws.Range("H24").Value = -59.82809513920385

$ws.Range("G25").Value = -0.2095999551699858
$ws.Range("H25").Value = 5.778330145632152

$ws.Range("G26").Value = 0.1582152947984669
$ws.Range("H26").Value = -0.4821450170865214

$ws.Range("G27").Value = 0.2213303484106095
$ws.Range("H27").Value = 10.4077384612342

$ws.Range("G28").Value = 0.01760782440066697
$ws.Range("H28").Value = 319.2632285717464

$ws.Range("G29").Value = 0.01891004639548182
$ws.Range("H29").Value = 22.97224141966725
